$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: I1 = "I0", J1 = "IF" ---
# Copy the formatting of the existing "IP" header (H1, style index 1:
# bold font, thin border all around, centered/top aligned) onto the two
# new header cells before writing their text, so I1/J1 end up styled the
# same way as the other header cells (B1..H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-25: I = 1 (constant), J = same value as H (IP) ---
for ($r = 2; $r -le 25; $r++) {
    $ip = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ip
}
